$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.944.99"
$ws.Range("E2").Value = "  +0.19%  "

$ws.Range("D3").Value = "2.264.99"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("E4").Value = "  -0.73%  "

$ws.Range("D5").Value = "'301.52"
$ws.Range("E5").Value = "  -2.00%  "

$ws.Range("D6").Value = "'93.79"
$ws.Range("E6").Value = "  -2.57%  "

$ws.Range("E7").Value = "  -1.69%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.61%  "

$ws.Range("E9").Value = "  -3.08%  "

$ws.Range("D10").Value = "'34.08"
$ws.Range("E10").Value = "  -4.09%  "

$ws.Range("D11").Value = "'0.0787"
$ws.Range("E11").Value = "  -2.58%  "

$ws.Range("D12").Value = "'7.19"
$ws.Range("E12").Value = "  -1.24%  "

$ws.Range("E13").Value = "  -1.18%  "

$ws.Range("D14").Value = "2.609.27"
$ws.Range("E14").Value = "  +0.51%  "

$ws.Range("D15").Value = "2.264.44"
$ws.Range("E15").Value = "  -0.77%  "

$ws.Range("E16").Value = "  +0.22%  "

$ws.Range("D17").Value = "'0.796"
$ws.Range("E17").Value = "  -5.61%  "

$ws.Range("D18").Value = "44.758.25"
$ws.Range("E18").Value = "  +0.44%  "

$ws.Range("D19").Value = "'12.98"
$ws.Range("E19").Value = "  +7.64%  "

$ws.Range("D20").Value = "0.0₃0920"
$ws.Range("E20").Value = "  -3.74%  "

$ws.Range("E21").Value = "  -3.90%  "

$ws.Range("D22").Value = "'65.37"
$ws.Range("E22").Value = "  -0.42%  "

$ws.Range("D23").Value = "'238.28"
$ws.Range("E23").Value = "  -0.31%  "

$ws.Range("E24").Value = "  -3.35%  "

$ws.Range("D25").Value = "'0.997"
$ws.Range("E25").Value = "  -0.50%  "

$ws.Range("D26").Value = "'1.90"
$ws.Range("E26").Value = "  -5.21%  "

$ws.Range("D27").Value = "'41.00"
$ws.Range("E27").Value = "  +9.66%  "

$ws.Range("D28").Value = "'2.30"
$ws.Range("E28").Value = "  +0.87%  "

$ws.Range("D29").Value = "'9.53"
$ws.Range("E29").Value = "  -3.89%  "

$ws.Range("D30").Value = "'19.49"
$ws.Range("E30").Value = "  -2.61%  "

$ws.Range("D31").Value = "'151.63"

$ws.Range("E32").Value = "  -8.03%  "

$ws.Range("D33").Value = "'0.0787"
$ws.Range("E33").Value = "  -1.88%  "

$ws.Range("D34").Value = "'2.55"
$ws.Range("E34").Value = "  -3.06%  "

$ws.Range("E35").Value = "  -6.47%  "

$ws.Range("E36").Value = "  -2.17%  "

$ws.Range("D37").Value = "'0.104"
$ws.Range("E37").Value = "  -5.46%  "

$ws.Range("D38").Value = "'1.75"
$ws.Range("E38").Value = "  -6.29%  "

$ws.Range("D39").Value = "'0.0313"
$ws.Range("E39").Value = "  +2.36%  "

$ws.Range("D40").Value = "'3.77"
$ws.Range("E40").Value = "  -1.77%  "

$ws.Range("D41").Value = "'3.24"
$ws.Range("E41").Value = "  -5.39%  "

$ws.Range("D42").Value = "'13.59"
$ws.Range("E42").Value = "  -9.26%  "

$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.93%  "

$ws.Range("E44").Value = "  +10.44%  "

$ws.Range("D45").Value = "1.764.90"
$ws.Range("E45").Value = "  -4.29%  "

$ws.Range("E46").Value = "  +1.05%  "

$ws.Range("D47").Value = "'70.35"
$ws.Range("E47").Value = "  -0.45%  "

$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "'75.36"
$ws.Range("E48").Value = "  -5.91%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'96.52"
$ws.Range("E49").Value = "  -3.20%  "

$ws.Range("D50").Value = "2.488.22"
$ws.Range("E50").Value = "  +0.50%  "

$ws.Range("D51").Value = "'52.93"
$ws.Range("E51").Value = "  -4.12%  "

